# Applies the updated leve profit figures across all job sheets
# as produced by the scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1749.5
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 2499
$ws.Range("K51").Value = 1000
$ws.Range("L51").Value = 2499
$ws.Range("M51").Value = -516
$ws.Range("N51").Value = -3467

$ws.Range("H129").Value = 2133.6155
$ws.Range("I129").Value = 1612.5454
$ws.Range("K129").Value = 4837.6362
$ws.Range("M129").Value = 162.3638000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1831.421
$ws.Range("I2").Value = 1736.3793
$ws.Range("K2").Value = 1736.3793
$ws.Range("M2").Value = -1623.3793

$ws.Range("H32").Value = 5623196.5
$ws.Range("I32").Value = 5819209.5
$ws.Range("J32").Value = 4166.6665
$ws.Range("K32").Value = 5819209.5
$ws.Range("L32").Value = 4166.6665
$ws.Range("M32").Value = -5818922.5
$ws.Range("N32").Value = -4740.6665

$ws.Range("H61").Value = 766189.2
$ws.Range("I61").Value = 1048354
$ws.Range("K61").Value = 1048354
$ws.Range("M61").Value = -1048142

$ws.Range("H74").Value = 2051838.9
$ws.Range("J74").Value = 3907.7
$ws.Range("L74").Value = 3907.7
$ws.Range("N74").Value = -5655.7

$ws.Range("H77").Value = 2051838.9
$ws.Range("J77").Value = 3907.7
$ws.Range("L77").Value = 19538.5
$ws.Range("N77").Value = -28274.5

$ws.Range("H97").Value = 696.6667
$ws.Range("I97").Value = 696.6667
$ws.Range("K97").Value = 696.6667
$ws.Range("M97").Value = -200.6667

$ws.Range("H110").Value = 2184
$ws.Range("I110").Value = 2218.2
$ws.Range("K110").Value = 2218.2
$ws.Range("M110").Value = -173.1999999999998

$ws.Range("H116").Value = 1831.421
$ws.Range("I116").Value = 1736.3793
$ws.Range("K116").Value = 1736.3793
$ws.Range("M116").Value = 557.6206999999999

$ws.Range("H132").Value = 268368
$ws.Range("I132").Value = 417765.34
$ws.Range("J132").Value = 6922.6875
$ws.Range("K132").Value = 1253296.02
$ws.Range("L132").Value = 20768.0625
$ws.Range("M132").Value = -1250766.02
$ws.Range("N132").Value = -25828.0625

$ws.Range("H136").Value = 766189.2
$ws.Range("I136").Value = 1048354
$ws.Range("K136").Value = 3145062
$ws.Range("M136").Value = -3142512

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1831.421
$ws.Range("I3").Value = 1736.3793
$ws.Range("K3").Value = 1736.3793
$ws.Range("M3").Value = -1622.3793

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9762.625
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 9762.625
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 9762.625
$ws.Range("N31").Value = -10352.625
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 9762.625
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9762.625
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 9762.625
$ws.Range("N34").Value = -10166.625
$ws.Range("M34").ClearContents()

$ws.Range("H43").Value = 100547.5
$ws.Range("J43").Value = 100547.5
$ws.Range("L43").Value = 100547.5
$ws.Range("N43").Value = -100915.5

$ws.Range("H74").Value = 63310.832
$ws.Range("J74").Value = 89995
$ws.Range("L74").Value = 89995
$ws.Range("N74").Value = -91743

$ws.Range("H77").Value = 63310.832
$ws.Range("J77").Value = 89995
$ws.Range("L77").Value = 269985
$ws.Range("N77").Value = -278721

$ws.Range("H87").Value = 123553.336
$ws.Range("J87").Value = 123553.336
$ws.Range("L87").Value = 123553.336
$ws.Range("N87").Value = -125925.336

$ws.Range("H90").Value = 123553.336
$ws.Range("J90").Value = 123553.336
$ws.Range("L90").Value = 370660.008
$ws.Range("N90").Value = -382516.008

$ws.Range("H101").Value = 100547.5
$ws.Range("J101").Value = 100547.5
$ws.Range("L101").Value = 100547.5
$ws.Range("N101").Value = -107037.5

$ws.Range("H105").Value = 1278.6
$ws.Range("I105").Value = 1278.6
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1278.6
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 468.4000000000001
$ws.Range("N105").ClearContents()

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H135").Value = 97243.75
$ws.Range("J135").Value = 97243.75
$ws.Range("L135").Value = 97243.75
$ws.Range("N135").Value = -107383.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5942.154
$ws.Range("I3").Value = 1632
$ws.Range("K3").Value = 4896
$ws.Range("M3").Value = -4784

$ws.Range("H33").Value = 563
$ws.Range("I33").Value = 434
$ws.Range("K33").Value = 2604
$ws.Range("M33").Value = -2321

$ws.Range("H134").Value = 1529.625
$ws.Range("I134").Value = 1529.625
$ws.Range("K134").Value = 4588.875
$ws.Range("M134").Value = 481.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2358.7144
$ws.Range("I122").Value = 1603.0605
$ws.Range("K122").Value = 4809.181500000001
$ws.Range("M122").Value = -2359.181500000001

$ws.Range("H132").Value = 229050.7
$ws.Range("I132").Value = 356024.12
$ws.Range("J132").Value = 1835.1052
$ws.Range("K132").Value = 1068072.36
$ws.Range("L132").Value = 5505.3156
$ws.Range("M132").Value = -1065542.36
$ws.Range("N132").Value = -10565.3156

$ws.Range("H139").Value = 138181.64
$ws.Range("J139").Value = 145999.8
$ws.Range("L139").Value = 145999.8
$ws.Range("N139").Value = -156279.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1911
$ws.Range("I82").Value = 1776.3
$ws.Range("J82").Value = 2079.375
$ws.Range("K82").Value = 1776.3
$ws.Range("L82").Value = 2079.375
$ws.Range("M82").Value = -1415.3
$ws.Range("N82").Value = -2801.375

$ws.Range("H85").Value = 1911
$ws.Range("I85").Value = 1776.3
$ws.Range("J85").Value = 2079.375
$ws.Range("K85").Value = 1776.3
$ws.Range("L85").Value = 2079.375
$ws.Range("M85").Value = -528.3
$ws.Range("N85").Value = -4575.375

$ws.Range("H122").Value = 40533.594
$ws.Range("I122").Value = 3183.0476
$ws.Range("J122").Value = 171260.5
$ws.Range("K122").Value = 9549.1428
$ws.Range("L122").Value = 513781.5
$ws.Range("M122").Value = -7099.1428
$ws.Range("N122").Value = -518681.5

$ws.Range("H132").Value = 542703.7
$ws.Range("I132").Value = 667208.5
$ws.Range("J132").Value = 3182.75
$ws.Range("K132").Value = 2001625.5
$ws.Range("L132").Value = 9548.25
$ws.Range("M132").Value = -1999095.5
$ws.Range("N132").Value = -14608.25

$ws.Range("H135").Value = 48333
$ws.Range("J135").Value = 48333
$ws.Range("L135").Value = 48333
$ws.Range("N135").Value = -58473

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2707.818
$ws.Range("I122").Value = 2279.1738
$ws.Range("J122").Value = 3693.7
$ws.Range("K122").Value = 6837.5214
$ws.Range("L122").Value = 11081.1
$ws.Range("M122").Value = -4387.5214
$ws.Range("N122").Value = -15981.1

$ws.Range("H123").Value = 47000
$ws.Range("J123").Value = 47000
$ws.Range("L123").Value = 47000
$ws.Range("N123").Value = -56800

$ws.Range("H132").Value = 3356712.2
$ws.Range("I132").Value = 4376922.5
$ws.Range("K132").Value = 13130767.5
$ws.Range("M132").Value = -13128237.5
